$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the dates in column F (rows 2-7) forward by 21 days (3 weeks),
# preserving the existing date number format / style on those cells.
$ws.Range("F2").Value = 44914
$ws.Range("F3").Value = 44913
$ws.Range("F4").Value = 44912
$ws.Range("F5").Value = 44911
$ws.Range("F6").Value = 44910
$ws.Range("F7").Value = 44909
